$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Insert a new bullet "Improved commenting of components." (ilvl=1)
#    right before the "Unit tests." bullet, and move the _GoBack
#    bookmark so it ends up right after the new bullet's text.
# -----------------------------------------------------------------

# Locate the "Unit tests." paragraph (short paragraph so we don't
# accidentally match the longer one that also contains "unit tests").
$unitTestsIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($t -like "*Unit tests*" -and $t.Length -lt 20) {
        $unitTestsIndex = $i
        break
    }
}

$targetRange = $d.Paragraphs.Item($unitTestsIndex).Range
$insertPoint = $d.Range($targetRange.Start, $targetRange.Start)
# A trailing placeholder character ("X") is inserted along with the
# new text/paragraph mark so that the paragraph end position is never
# used on its own as a zero-width anchor (avoids a boundary case in
# the range engine). It is stripped again once the bookmark has been
# anchored around it, collapsing the bookmark to the correct spot.
$insertPoint.InsertBefore("Improved commenting of components.X`r")

$newPara = $d.Paragraphs.Item($unitTestsIndex)
$newPara.Range.ListFormat.ListLevelNumber = 2

$placeholderStart = $newPara.Range.End - 2
$placeholderRange = $d.Range($placeholderStart, $placeholderStart + 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange)
$placeholderRange = $d.Range($placeholderStart, $placeholderStart + 1)
$placeholderRange.Delete()

# -----------------------------------------------------------------
# 2. Move the lastRenderedPageBreak marker from the start of the
#    "Important! Ensure -fno-strict-aliasing..." run to the start of
#    the "Python bindings." run (the immediately preceding bullet).
# -----------------------------------------------------------------

$pyIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($t -like "*Python bindings*") {
        $pyIndex = $i
        break
    }
}
$pyPara = $d.Paragraphs.Item($pyIndex)
$pyRange = $d.Range($pyPara.Range.Start, $pyPara.Range.End - 1)
$pyRange.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Python bindings.</w:t></w:r></w:p>")

$aliasIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($t -like "*fno-strict-aliasing*") {
        $aliasIndex = $i
        break
    }
}
$aliasPara = $d.Paragraphs.Item($aliasIndex)
$aliasRange = $d.Range($aliasPara.Range.Start, $aliasPara.Range.End - 1)
$aliasRange.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t>Important! Ensure -fno-strict-aliasing is used under GCC as it seems Boost.Python has aliasing violations which cause spurious segfaults and other issues.</w:t></w:r></w:p>")
